$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: new entry - 2021-09-06, 2 hours, "Meeting"
# (A27 already carries the date-formatted style in the template, so a plain
# value assignment keeps it.)
$ws.Range("A27").Value = 44445
$ws.Range("B27").Value = 2
$ws.Range("D27").Value = "Meeting"

# Row 28: new entry - 2021-09-08, 6 hours, "Chart navigation (mobile support)"
$ws.Range("A28").Value = 44447
$ws.Range("B28").Value = 6
$ws.Range("D28").Value = "Chart navigation (mobile support)"

# Row 29: new entry - 2021-09-09, 1 hour, "Meeting"
# A29 didn't previously exist as a styled cell, so pull the date style down
# from the row above first (copy) and then overwrite with the real value.
$ws.Range("A28").Copy($ws.Range("A29"))
$ws.Range("A29").Value = 44448
$ws.Range("B29").Value = 1
$ws.Range("D29").Value = "Meeting"

# Recalculate so the running-total (Zwischensumme) formulas in column C pick
# up the new hours entered above (rows 27-32 all chain off of each other).
$excel.Calculate()

# Move the active selection from L7 to H16.
$ws.Range("H16").Select()
